$d = $word.ActiveDocument

function New-PkgXml([string]$bodyInner) {
    return '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Set-ParaXml([int]$index, [string]$pInner) {
    $wrapped = '<w:p>' + $pInner + '</w:p>'
    $full = New-PkgXml $wrapped
    $rng = $d.Paragraphs($index).Range
    $rng.InsertXML($full)
}

# --- "scipy" / "numpy" / "matplotlib" (paragraphs 6-8): wrap in spellStart/gramStart .. spellEnd/gramEnd ---
$t6 = '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>scipy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>'
Set-ParaXml 6 $t6

$t7 = '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>numpy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>'
Set-ParaXml 7 $t7

$t8 = '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>matplotlib</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>'
Set-ParaXml 8 $t8

# --- "1. correlation matrix" (paragraph 10) ---
$t10 = (
    '<w:r><w:t xml:space="preserve">1. </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>corre</w:t></w:r>' +
    '<w:r><w:t>l</w:t></w:r>' +
    '<w:r><w:t>a</w:t></w:r>' +
    '<w:r><w:t>tion</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> matrix</w:t></w:r>'
)
Set-ParaXml 10 $t10

# --- "2. for each opponent, find wins, loses and ties" (paragraph 11) ---
$t11 = (
    '<w:r><w:t xml:space="preserve">2. </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>for</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> each opponent</w:t></w:r>' +
    '<w:r><w:t>,</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> find wins, loses and ties</w:t></w:r>'
)
Set-ParaXml 11 $t11

# --- "3. what was the magnitude of win, loss or tie" (paragraph 12) ---
$t12 = (
    '<w:r><w:t xml:space="preserve">3. </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>what</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> was the magnitude of win, loss or tie</w:t></w:r>'
)
Set-ParaXml 12 $t12

# --- "4. when i)win look at conversion, pass diff, poss diff, ruck maul diff etc.." (paragraph 13) ---
$t13 = (
    '<w:r><w:t xml:space="preserve">4. </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>when</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>i</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">)win look at conversion, pass diff, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>poss</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> diff, ruck maul diff etc..</w:t></w:r>'
)
Set-ParaXml 13 $t13

# --- "<tab>ii) loss, same..." (paragraph 14) ---
$t14 = (
    '<w:r><w:tab/></w:r>' +
    '<w:r><w:t xml:space="preserve">ii) </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>loss</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>, same&#8230;</w:t></w:r>'
)
Set-ParaXml 14 $t14

# --- "<tab>iii) tie, same" + bookmark (paragraph 15) ---
$t15 = (
    '<w:r><w:tab/></w:r>' +
    '<w:r><w:t xml:space="preserve">iii) </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>tie</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>, same</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>'
)
Set-ParaXml 15 $t15

# --- "5. Look at for e.x against Fiji, lost most matches" (paragraph 16) ---
$t16 = (
    '<w:r><w:t xml:space="preserve">5. Look at for </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>e.x</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> against Fiji, lost most matches</w:t></w:r>'
)
Set-ParaXml 16 $t16

# --- "6. which other team won most matches against fiji" (paragraph 17) ---
$t17 = (
    '<w:r><w:t xml:space="preserve">6. </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>which</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> other team won most matches against </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>fiji</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
)
Set-ParaXml 17 $t17

# --- "<tab>7. repat number 4 I) ii) and iii)" (paragraph 18) ---
$t18 = (
    '<w:r><w:tab/></w:r>' +
    '<w:r><w:t xml:space="preserve">7. </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>repat</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> number 4 I) ii) and iii)</w:t></w:r>'
)
Set-ParaXml 18 $t18

# --- "*repeat" (paragraph 19) ---
$t19 = (
    '<w:r><w:t>*</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>repeat</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>'
)
Set-ParaXml 19 $t19

# --- Replace the middle blank paragraph (23) with 3 new paragraphs:
#     "Restarts - take median, median of wins" (red)
#     "Hypothesis is > mean = win" (red)
#     empty paragraph with green run-properties mark
$newParas = (
    '<w:p><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Restarts &#8211; take median, median of wins</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Hypothesis is &gt; mean = win</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:rPr><w:color w:val="008000"/></w:rPr></w:pPr></w:p>'
)
$full23 = New-PkgXml $newParas
$p23 = $d.Paragraphs(23).Range
$p23.InsertXML($full23)

Write-Output "done"
